# Update with Correct Forecast output
#
# - Renames Sheet1 to "Sales vs PO".
# - Inserts a new "Order Week" column (C) on that sheet: the old "ds" values
#   move into "Order Week", "ds" itself shifts forward 6 days (the
#   order-placed date becomes roughly a week after the order week),
#   and the old "PO_Requested_Qty" column becomes column D, zeroed out.
# - Adds three brand-new sheets: "Weekly Growth" (per-week PO qty + WoW
#   growth %), "Volume Insights" (aggregate PO stats) and "Prediction Info"
#   (next week's forecasted PO quantity).

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Sales vs PO" (was "Sheet1")
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# Remember the original header/date formatting before we start overwriting
# cells, so the brand new column (D) and new sheets can reuse them exactly.
$headerFormatCell = $ws1.Range("A1")
$dateFormatCell = $ws1.Range("A2")

# New header row: ds, y, Order Week, PO_Requested_Qty
$ws1.Cells.Item(1, 3).Value = "Order Week"
$ws1.Cells.Item(1, 4).Value = "PO_Requested_Qty"

# Make sure every header cell (including the brand-new D1) has the bold /
# bordered / centered header style used by the rest of the row.
$headerFormatCell.Copy()
$ws1.Range("A1:D1").PasteSpecial($xlPasteFormats)

# Data rows: ds shifts forward 6 days, Order Week takes the old ds value,
# and PO_Requested_Qty (col D) is reset to 0 for every row.
$sheet1Rows = @(
    @(45494, 0,   45488, 0),
    @(45501, 0,   45495, 0),
    @(45508, 0,   45502, 0),
    @(45515, 3,   45509, 0),
    @(45522, 18,  45516, 0),
    @(45529, 143, 45523, 0),
    @(45536, 28,  45530, 0),
    @(45543, 41,  45537, 0),
    @(45550, 19,  45544, 0),
    @(45557, 11,  45551, 0),
    @(45564, 22,  45558, 0),
    @(45571, 43,  45565, 0),
    @(45578, 27,  45572, 0),
    @(45585, 23,  45579, 0),
    @(45592, 28,  45586, 0),
    @(45599, 30,  45593, 0),
    @(45606, 16,  45600, 0),
    @(45613, 86,  45607, 0),
    @(45620, 258, 45614, 0),
    @(45627, 301, 45621, 0),
    @(45634, 84,  45628, 0),
    @(45641, 12,  45635, 0),
    @(45648, 11,  45642, 0),
    @(45655, 8,   45649, 0)
)

$r = 2
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Re-apply the yyyy-mm-dd date style to the "ds" and "Order Week" columns.
$dateFormatCell.Copy()
$ws1.Range("A2:A25").PasteSpecial($xlPasteFormats)
$dateFormatCell.Copy()
$ws1.Range("C2:C25").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# Sheet 2: "Weekly Growth"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$ws2.Cells.Item(1, 1).Value = "ds"
$ws2.Cells.Item(1, 2).Value = "PO_Requested_Qty"
$ws2.Cells.Item(1, 3).Value = "Growth%"

$headerFormatCell.Copy()
$ws2.Range("A1:C1").PasteSpecial($xlPasteFormats)

$sheet2Rows = @(
    @(45495, 16,  0),
    @(45502, 192, 1100),
    @(45509, 112, -41.66666666666666),
    @(45516, 16,  -85.71428571428572),
    @(45523, 32,  100),
    @(45530, 128, 300),
    @(45537, 464, 262.5),
    @(45544, 256, -44.82758620689656),
    @(45551, 336, 31.25),
    @(45579, 640, 90.47619047619047),
    @(45593, 16,  -97.5)
)

$r = 2
foreach ($row in $sheet2Rows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $r++
}

$dateFormatCell.Copy()
$ws2.Range("A2:A12").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# Sheet 3: "Volume Insights"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$ws3.Cells.Item(1, 1).Value = "Total_PO_Quantity"
$ws3.Cells.Item(1, 2).Value = "Average_PO_Quantity"
$ws3.Cells.Item(1, 3).Value = "Max_PO_Quantity"
$ws3.Cells.Item(1, 4).Value = "Min_PO_Quantity"

$headerFormatCell.Copy()
$ws3.Range("A1:D1").PasteSpecial($xlPasteFormats)

$ws3.Cells.Item(2, 1).Value = 2208
$ws3.Cells.Item(2, 2).Value = 200.7272727272727
$ws3.Cells.Item(2, 3).Value = 640
$ws3.Cells.Item(2, 4).Value = 16

# ---------------------------------------------------------------------------
# Sheet 4: "Prediction Info"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

$ws4.Cells.Item(1, 1).Value = "Predicted_Next_Week_PO_Quantity"
$headerFormatCell.Copy()
$ws4.Range("A1").PasteSpecial($xlPasteFormats)

$ws4.Cells.Item(2, 1).Value = 384.8727272727274

# ---------------------------------------------------------------------------
# Re-select the first sheet / first cell to mirror the original view state.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A1").Select()
